$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3999.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3999.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3999.5
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4967.5

$ws.Range("H58").Value = 1588
$ws.Range("I58").Value = 61.714287
$ws.Range("J58").Value = 3114.2856
$ws.Range("K58").Value = 185.142861
$ws.Range("L58").Value = 9342.856800000001
$ws.Range("M58").Value = -35.14286099999998
$ws.Range("N58").Value = -9642.856800000001

$ws.Range("H116").Value = 5872.5
$ws.Range("I116").Value = 5548.2856
$ws.Range("K116").Value = 5548.2856
$ws.Range("M116").Value = -2106.2856

$ws.Range("H138").Value = 2156.5833
$ws.Range("I138").Value = 1409.3334
$ws.Range("K138").Value = 4228.0002
$ws.Range("M138").Value = 911.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 188.2
$ws.Range("I5").Value = 146.33333
$ws.Range("J5").Value = 251
$ws.Range("K5").Value = 146.33333
$ws.Range("L5").Value = 251
$ws.Range("M5").Value = -34.33332999999999
$ws.Range("N5").Value = -475

$ws.Range("H37").Value = 4684.6665
$ws.Range("I37").Value = 4684.6665
$ws.Range("K37").Value = 4684.6665
$ws.Range("M37").Value = -4411.6665

$ws.Range("H74").Value = 2286.8928
$ws.Range("I74").Value = 2040.1923
$ws.Range("J74").Value = 5494
$ws.Range("K74").Value = 2040.1923
$ws.Range("L74").Value = 5494
$ws.Range("M74").Value = -1166.1923
$ws.Range("N74").Value = -7242

$ws.Range("H77").Value = 2286.8928
$ws.Range("I77").Value = 2040.1923
$ws.Range("J77").Value = 5494
$ws.Range("K77").Value = 10200.9615
$ws.Range("L77").Value = 27470
$ws.Range("M77").Value = -5832.961499999999
$ws.Range("N77").Value = -36206

$ws.Range("H122").Value = 4464.143
$ws.Range("I122").Value = 4124.75
$ws.Range("K122").Value = 12374.25
$ws.Range("M122").Value = -9924.25

$ws.Range("H132").Value = 807.8333
$ws.Range("I132").Value = 807.8333
$ws.Range("K132").Value = 2423.4999
$ws.Range("M132").Value = 106.5001000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 188.2
$ws.Range("I4").Value = 146.33333
$ws.Range("J4").Value = 251
$ws.Range("K4").Value = 146.33333
$ws.Range("L4").Value = 251
$ws.Range("M4").Value = -31.33332999999999
$ws.Range("N4").Value = -481

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 125.666664
$ws.Range("I7").Value = 83
$ws.Range("K7").Value = 83
$ws.Range("M7").Value = 30

$ws.Range("H31").Value = 5003.256
$ws.Range("I31").Value = 4397.027
$ws.Range("J31").Value = 8741.666999999999
$ws.Range("K31").Value = 4397.027
$ws.Range("L31").Value = 8741.666999999999
$ws.Range("M31").Value = -4102.027
$ws.Range("N31").Value = -9331.666999999999

$ws.Range("H34").Value = 5003.256
$ws.Range("I34").Value = 4397.027
$ws.Range("J34").Value = 8741.666999999999
$ws.Range("K34").Value = 4397.027
$ws.Range("L34").Value = 8741.666999999999
$ws.Range("M34").Value = -4195.027
$ws.Range("N34").Value = -9145.666999999999

$ws.Range("H41").Value = 49505
$ws.Range("I41").Value = 245
$ws.Range("J41").Value = 98765
$ws.Range("K41").Value = 245
$ws.Range("L41").Value = 98765
$ws.Range("M41").Value = 183
$ws.Range("N41").Value = -99621

$ws.Range("H74").Value = 76153
$ws.Range("J74").Value = 76153
$ws.Range("L74").Value = 76153
$ws.Range("N74").Value = -77901

$ws.Range("H77").Value = 76153
$ws.Range("J77").Value = 76153
$ws.Range("L77").Value = 228459
$ws.Range("N77").Value = -237195

$ws.Range("H99").Value = 5805.4287
$ws.Range("I99").Value = 5129.8
$ws.Range("K99").Value = 5129.8
$ws.Range("M99").Value = -3631.8

$ws.Range("H126").Value = 5805.4287
$ws.Range("I126").Value = 5129.8
$ws.Range("K126").Value = 15389.4
$ws.Range("M126").Value = -12919.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 350.15384
$ws.Range("I12").Value = 25.666666
$ws.Range("J12").Value = 447.5
$ws.Range("K12").Value = 76.99999800000001
$ws.Range("L12").Value = 1342.5
$ws.Range("M12").Value = 96.00000199999999
$ws.Range("N12").Value = -1688.5

$ws.Range("H140").Value = 2475.3635
$ws.Range("I140").Value = 1814.3334
$ws.Range("K140").Value = 5443.0002
$ws.Range("M140").Value = -263.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 52499.5
$ws.Range("J93").Value = 52499.5
$ws.Range("L93").Value = 52499.5
$ws.Range("N93").Value = -56243.5

$ws.Range("H132").Value = 50406.25
$ws.Range("J132").Value = 9832.333000000001
$ws.Range("L132").Value = 29496.999
$ws.Range("N132").Value = -34556.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 716.6667
$ws.Range("I16").Value = 716.6667
$ws.Range("K16").Value = 716.6667
$ws.Range("M16").Value = -546.6667

$ws.Range("H22").Value = 2333.3333
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2333.3333
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2333.3333
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2923.3333

$ws.Range("H27").Value = 2333.3333
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2333.3333
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2333.3333
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2547.3333

$ws.Range("H40").Value = 8171.8423
$ws.Range("I40").Value = 7285.467
$ws.Range("K40").Value = 7285.467
$ws.Range("M40").Value = -7149.467

$ws.Range("H82").Value = 5251.5
$ws.Range("I82").Value = 4305
$ws.Range("J82").Value = 6829
$ws.Range("K82").Value = 4305
$ws.Range("L82").Value = 6829
$ws.Range("M82").Value = -3944
$ws.Range("N82").Value = -7551

$ws.Range("H85").Value = 5251.5
$ws.Range("I85").Value = 4305
$ws.Range("J85").Value = 6829
$ws.Range("K85").Value = 4305
$ws.Range("L85").Value = 6829
$ws.Range("M85").Value = -3057
$ws.Range("N85").Value = -9325

$ws.Range("H136").Value = 6785.3335
$ws.Range("I136").Value = 6383.5
$ws.Range("K136").Value = 19150.5
$ws.Range("M136").Value = -16600.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 733
$ws.Range("I100").Value = 599.5
$ws.Range("K100").Value = 1199
$ws.Range("M100").Value = -658
